$d = $word.ActiveDocument

# --- Part 1: capitalize "salut" -> "Salut", keeping "S" and "alut" as two
#     separate runs (matches the diff's two <w:r> elements) ---

# Replace the lower-case "s" with an upper-case "S" (still a single run at
# this point).
$d.Range(0, 1).Text = "S"

# Force a run boundary right after the "S" by splitting the paragraph there
# and then deleting the freshly inserted paragraph mark again, which merges
# the two paragraphs back into one while leaving the two runs distinct.
$split = $d.Range(1, 1)
$split.InsertParagraphAfter()
$firstPara = $d.Paragraphs(1).Range
$newMark = $d.Range($firstPara.End - 1, $firstPara.End)
$newMark.Delete()

# --- Part 2: add a new paragraph "Projet devops" right after "Salut" ---
# Locate the end of the (now merged) first paragraph's text, i.e. right
# before the "_GoBack" bookmark, and insert a paragraph break followed by
# the new line's text there. Using a literal CR makes the break carry the
# trailing bookmark into the newly created second paragraph, exactly like
# the target document.
$endOfFirstPara = $d.Paragraphs(1).Range.End - 1
$breakPoint = $d.Range($endOfFirstPara, $endOfFirstPara)
$breakPoint.InsertBefore([char]13)

$newLineStart = $endOfFirstPara + 1
$d.Range($newLineStart, $newLineStart).InsertBefore("Projet devops")
